$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 102
$ws.Range("H4").Value = 1135
$ws.Range("I4").Value = 1001
$ws.Range("J4").Value = 1108
$ws.Range("Q4").Value = 718
